# Update division problems to new values as per commit "Update master to
# output generated at 9a8706d".
$d = $word.ActiveDocument

$replacements = @(
    @("271÷7=", "202÷4="),
    @("796÷7=", "509÷2="),
    @("112÷4=", "695÷7="),
    @("922÷3=", "269÷5="),
    @("609÷5=", "130÷2="),
    @("778÷8=", "297÷8="),
    @("231÷7=", "721÷2="),
    @("814÷7=", "582÷4="),
    @("182÷3=", "294÷4="),
    @("440÷6=", "218÷7="),
    @("734÷6=", "580÷2="),
    @("232÷7=", "678÷6="),
    @("110÷5=", "341÷9="),
    @("764÷4=", "936÷5="),
    @("419÷7=", "572÷6="),
    @("128÷9=", "734÷2="),
    @("904÷4=", "993÷2="),
    @("587÷7=", "523÷7="),
    @("102÷8=", "723÷9="),
    @("327÷2=", "789÷5="),
    @("466÷6=", "828÷2="),
    @("292÷9=", "487÷9="),
    @("314÷2=", "989÷7="),
    @("104÷9=", "975÷9="),
    @("973÷6=", "312÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
